$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Add New Bank Accounts", "PASSED", "chrome"),
    @("Edit The Bank Accounts", "PASSED", "chrome"),
    @("Delete The Bank Accounts", "PASSED", "chrome"),
    @("Add New Bank Accounts", "PASSED", "chrome"),
    @("Edit The Bank Accounts", "PASSED", "chrome"),
    @("Delete The Bank Accounts", "PASSED", "chrome"),
    @("", "FAILED", "chrome"),
    @("", "PASSED", "chrome")
)

$row = 14
foreach ($r in $data) {
    $colA = $ws.Cells.Item($row, 1)
    if ($r[0] -eq "") {
        # Writing a literal empty string via .Value clears the cell instead
        # of leaving an empty-string text cell behind, so force a text
        # value through the quote-prefix trick, then strip the resulting
        # quote-prefix style back off so the cell looks like a plain
        # empty shared-string cell again.
        $colA.Value = "'"
        $colA.Style = "Normal"
    } else {
        $colA.Value = $r[0]
    }
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
